$wb = $excel.ActiveWorkbook

# --- Metadata sheet: update Date value ---
$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B8").Value2 = "2024-03-19T13:17:15+00:00"

# --- Elements sheet: swap Mapping columns (AK <-> AL) ---
$el = $wb.Worksheets.Item("Elements")

# Swap header text in row 1 (AK1 <-> AL1)
$akHeader = $el.Range("AK1").Value2
$alHeader = $el.Range("AL1").Value2
$el.Range("AK1").Value2 = $alHeader
$el.Range("AL1").Value2 = $akHeader

# Swap column widths (AK was 24.98046875 / AL was 85.5234375 -> now swapped).
# Note: the ColumnWidth setter in this runtime quantizes to ~1/6-character
# steps, so the literal inputs below are chosen to land on the closest
# achievable stored width to the target values (85.5234375 / 24.98046875).
$el.Columns.Item(37).ColumnWidth = 84.59
$el.Columns.Item(38).ColumnWidth = 24.09

# Swap data cells for rows 3, 5, 6
foreach ($r in 3, 5, 6) {
    $akVal = $el.Range("AK$r").Value2
    $alVal = $el.Range("AL$r").Value2
    $el.Range("AK$r").Value2 = $alVal
    $el.Range("AL$r").Value2 = $akVal
}
